$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-word the existing "sportswriter" question in A22: the trailing
#    "This is an\nexample of" became "This is an example of" (line break
#    removed, joined onto one line). Everything else in that row (answer
#    choices, answer key) is left untouched.
$ws.Cells.Item(22, 1).Value = "A sportswriter wants to know how strongly residents support building a new stadium downtown for the local major league soccer team. She prints a survey in her column and asks her readers to send in their response. One thousand readers sent in their response.`nCritics of the poll argue that the poll only sampled readers`nof her column and not all residents of the city. This is an example of"

# 2) Append new "quiz 6" content as two new rows at the bottom of the sheet.

# Row 42: "It is difficult to establish the causal link ..." question.
$ws.Cells.Item(42, 1).Value = "It is difficult to establish the causal link between cigarette smoking and lung cancer because"
$ws.Cells.Item(42, 2).Value = "random allocation of subjects to smoking is unethical."
$ws.Cells.Item(42, 3).Value = "random allocation of subjects to smoking is unethical."
$ws.Cells.Item(42, 4).Value = "those who choose to smoke may be genetically at greater risk for lung cancer than those who don't choose to smoke."
$ws.Cells.Item(42, 5).Value = "experiments done on animals may not be valid for humans."
$ws.Cells.Item(42, 6).Value = "all of the above"

# Row 43: "A professor believes that students who smoke cigarettes ..." question.
$ws.Cells.Item(43, 1).Value = "A professor believes that students who smoke cigarettes tend to have lower grades. He collects data from 1326 randomly selected students at his university and discovers that, on average, students who smoke cigarettes do indeed tend to have lower grade point averages than students who do not smoke.`nThis study was based on"
$ws.Cells.Item(43, 2).Value = "a randomized comparative experiment."
$ws.Cells.Item(43, 3).Value = "a matched pairs experiment."
$ws.Cells.Item(43, 4).Value = "a voluntary response sample."
$ws.Cells.Item(43, 5).Value = "a probability sample."

# The long question-text cells (column A) use the wrapped-text style like the
# other question rows above them.
$ws.Range("A43").WrapText = $true

# Move the selection/scroll position to reflect where the author ended up
# after adding the new rows.
$ws.Range("G50").Select() | Out-Null
